$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.07%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.50%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.680"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-6.15%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08073"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.07%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.045"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.71%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.750"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.36%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.537"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.48%"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.91%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1945"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.57%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.779"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-16.47%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09382"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.20%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03749"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "6.97%"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.34%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001299"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.22%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006200"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.87%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.41%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.81%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.29%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2657"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.86%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04420"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.52%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001259"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.16%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004345"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.84%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001241"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "13.64%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02866"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "15.03%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05488"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.23%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007771"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.37%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009967"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.17%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1420"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.64%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002227"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.33%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01111"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.06%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006773"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.72%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.10%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002282"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "26.50%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003023"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-13.68%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.10%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.10%"
